# Update "想去人数" (want-to-go count) figures for several events.
# These values changed between scrapes (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13127
$ws1.Range("F7").Value = 58
$ws1.Range("F10").Value = 13083
$ws1.Range("F13").Value = 8783

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13127
$ws4.Range("F8").Value = 58
$ws4.Range("F11").Value = 13083
$ws4.Range("F14").Value = 8783
